$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1327811", "https://aiesec.org/opportunity/global-talent/1327811", "Software Engineering Intern", "Colombo, Sri Lanka", "No", "131 applicants", "3 - 6 Months", "Envision Circle (Pvt) Ltd"),
    @("1327658", "https://aiesec.org/opportunity/global-talent/1327658", "DevOps Engineer", "El-Kom El-Ahmar, Shibin el-Qanater, Al-Qalyubia Governorate, Egypt", "No", "18 applicants", "9 - 12 Weeks", "Etolv"),
    @("1324164", "https://aiesec.org/opportunity/global-talent/1324164", "Receptionist/ Hostess", "Nuwara Eliya, Sri Lanka", "No", "70 applicants", "3 - 6 Months", "Pedro Barn pvt ltd")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 1
        $cell = $ws.Cells.Item($row, $col)
        # Force text storage (e.g. numeric-looking opportunity IDs) then drop
        # the temporary number-format style so no extra formatting is left
        # behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $data[$i][$j]
        $cell.ClearFormats()
    }
}

# Excel's ColumnWidth -> stored OOXML "width" conversion adds a constant
# ~0.8333 character offset (5px at default Calibri 11 metrics). Subtract it
# here so the persisted <col width="..."/> exactly matches the target value.
$ws.Columns.Item(2).ColumnWidth = 54.166666666666664
$ws.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws.Columns.Item(4).ColumnWidth = 68.16666666666667
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666
$ws.Columns.Item(8).ColumnWidth = 27.166666666666668
